$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Ark1")
$ws2 = $wb.Worksheets.Item("Ark2")
$ws3 = $wb.Worksheets.Item("Ark3")

# Sheet1 (Ark1) data changes
$ws1.Range("E2").Value = 5
$ws1.Range("E3").Value = 5
$ws1.Range("E4").Value = 5
$ws1.Range("C5").Value = 864
$ws1.Range("E5").Value = 5
[void]$ws1.Range("E5").Select()

# Sheet2 (Ark2) - swap C2 and C3 values
$ws2.Range("C2").Value = "b"
$ws2.Range("C3").Value = "contract"

# Activate Ark2 so it becomes the tab-selected sheet, then select C3
[void]$ws2.Activate()
[void]$ws2.Range("C3").Select()

Write-Output "done"
